$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36. This shifts the existing rows 36-84
# down to 37-85 (standard Excel "insert row" behaviour, pushing cell
# content down while keeping row 35 and above untouched).
$ws.Rows("36:36").Insert()

# The newly inserted row 36 is blank; populate it with a duplicate of
# the record that is now sitting in row 37 (the data that used to be in
# row 36 before the insert), matching how the source row was cloned.
$ws.Range("A37:T37").Copy()
$ws.Range("A36").PasteSpecial()

# Finally, overwrite the two fields that differ in the new record.
$ws.Range("D36").Value = 44482
$ws.Range("M36").Value = 40
